# edit.ps1
# Applies the "Updated capital structure database" revision to the
# australia_computers_peripherals worksheet:
#   - Recomputed metrics for the existing Australia / Computers-Peripherals rows
#   - Row 2 company code changes from "4" to "3"
#   - Row 3 company is replaced by "Weebit Nano Limited (ASX:WBT)" and a number
#     of its margin columns (G,H,I,J,L,AO) are cleared out
#   - The separate old "Weebit Nano Limited" row (row 6) is removed entirely
#     since the dataset now only has 5 data rows (dimension A1:AQ6 -> A1:AQ5)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2 ---
$ws.Range("B2").NumberFormat = "@"   # keep as text, not a number
$ws.Range("B2").Value = "3"
$ws.Range("B2").Style = "Normal"     # drop the temporary text style
$ws.Range("G2").Value = -6.47244094488189
$ws.Range("H2").Value = -10.79396325459318
$ws.Range("I2").Value = -11.00393700787401
$ws.Range("J2").Value = -11.00393700787401
$ws.Range("K2").Value = -18.32
$ws.Range("L2").Value = -12.02099737532808
$ws.Range("U2").Value = 6.38
$ws.Range("V2").Value = 0.02274753093022426
$ws.Range("W2").Value = -1.234710743801653
$ws.Range("X2").Value = 0.06895119733028632
$ws.Range("Y2").Value = -1.303661941131939
$ws.Range("Z2").Value = 0.1333916849015317
$ws.Range("AA2").Value = -0.7137129109863674
$ws.Range("AB2").Value = 0.06758625691391859
$ws.Range("AC2").Value = -0.7814319310845015
$ws.Range("AD2").Value = 2.535
$ws.Range("AE2").Value = 0
$ws.Range("AF2").Value = 2.535
$ws.Range("AG2").Value = -3.845
$ws.Range("AH2").Value = 0.00895743891450681
$ws.Range("AI2").Value = 0.2712390327412797
$ws.Range("AJ2").Value = -0.0138996836873023
$ws.Range("AK2").Value = -1.296358732299393
$ws.Range("AL2").Value = 0.082
$ws.Range("AM2").Value = -0.074
$ws.Range("AN2").Value = -0.1541033434650456
$ws.Range("AO2").Value = -204.5121951219512
$ws.Range("AP2").Value = 0.2337386018237081
$ws.Range("AQ2").Value = 226.6216216216216

# --- Row 3 ---
$ws.Range("B3").NumberFormat = "@"   # keep as text, not a number
$ws.Range("B3").Value = "Weebit Nano Limited (ASX:WBT)"
$ws.Range("B3").Style = "Normal"     # drop the temporary text style
$ws.Range("K3").Value = -2.78
$ws.Range("U3").Value = 2.84
$ws.Range("V3").Value = 0.0129798903107861
$ws.Range("W3").Value = -2.482142857142857
$ws.Range("X3").Value = 0.06760457328546611
$ws.Range("Y3").Value = -2.549747430428323
$ws.Range("Z3").Value = -0
$ws.Range("AA3").Value = 57.80000000000021
$ws.Range("AB3").Value = 0.06758625691391859
$ws.Range("AC3").Value = 57.73241374308629
$ws.Range("AD3").Value = 0.08599999999999999
$ws.Range("AE3").Value = 0
$ws.Range("AF3").Value = 0.08599999999999999
$ws.Range("AG3").Value = -2.754
$ws.Range("AH3").Value = 0.0003928985864788063
$ws.Range("AI3").Value = 0.0358931552587646
$ws.Range("AJ3").Value = -0.0127472853003527
$ws.Range("AK3").Value = 6.202702702702704
$ws.Range("AL3").Value = 0
$ws.Range("AM3").Value = -0.113
$ws.Range("AN3").Value = -0.02986111111111111
$ws.Range("AP3").Value = 0.95625
$ws.Range("AQ3").Value = 25.57522123893805

# --- Row 4 ---
$ws.Range("G4").Value = -9.382716049382715
$ws.Range("H4").Value = -13.45679012345679
$ws.Range("I4").Value = -13.73456790123457
$ws.Range("J4").Value = -13.73456790123457
$ws.Range("K4").Value = -7.47
$ws.Range("L4").Value = -23.05555555555555
$ws.Range("U4").Value = 0.48
$ws.Range("V4").Value = 0.05177993527508091
$ws.Range("W4").Value = -1.234710743801653
$ws.Range("X4").Value = 0.07066518973660181
$ws.Range("Y4").Value = -1.305375933538255
$ws.Range("Z4").Value = 0.05196471531676023
$ws.Range("AA4").Value = -0.7137129109863674
$ws.Range("AB4").Value = 0.06771902009813409
$ws.Range("AC4").Value = -0.7814319310845015
$ws.Range("AD4").Value = 0.699
$ws.Range("AE4").Value = 0
$ws.Range("AF4").Value = 0.699
$ws.Range("AG4").Value = 0.219
$ws.Range("AH4").Value = 0.07011736382786639
$ws.Range("AI4").Value = 0.7681318681318682
$ws.Range("AJ4").Value = 0.023079355042681
$ws.Range("AK4").Value = 0.5093023255813953
$ws.Range("AL4").Value = 0.077
$ws.Range("AM4").Value = 0.077
$ws.Range("AN4").Value = -0.1603211009174312
$ws.Range("AO4").Value = -57.7922077922078
$ws.Range("AP4").Value = -0.05022935779816513
$ws.Range("AQ4").Value = -57.7922077922078

# --- Row 5 ---
$ws.Range("G5").Value = -3.375000000000001
$ws.Range("H5").Value = -7.675000000000001
$ws.Range("I5").Value = -7.858333333333333
$ws.Range("J5").Value = -7.858333333333333
$ws.Range("K5").Value = -8.07
$ws.Range("L5").Value = -6.725000000000001
$ws.Range("U5").Value = 3.06
$ws.Range("V5").Value = 0.0583969465648855
$ws.Range("W5").Value = -1.076
$ws.Range("X5").Value = 0.06895119733028632
$ws.Range("Y5").Value = -1.144951197330286
$ws.Range("Z5").Value = 0.2290076335877863
$ws.Range("AA5").Value = -1.799618320610687
$ws.Range("AB5").Value = 0.06740107553290864
$ws.Range("AC5").Value = -1.867019396143595
$ws.Range("AD5").Value = 1.75
$ws.Range("AE5").Value = 0
$ws.Range("AF5").Value = 1.75
$ws.Range("AG5").Value = -1.31
$ws.Range("AH5").Value = 0.03231763619575254
$ws.Range("AI5").Value = 0.2897350993377483
$ws.Range("AJ5").Value = -0.02564102564102564
$ws.Range("AK5").Value = -0.4395973154362416
$ws.Range("AL5").Value = 0.005
$ws.Range("AM5").Value = -0.038
$ws.Range("AN5").Value = -0.1900108577633007
$ws.Range("AO5").Value = -1886
$ws.Range("AP5").Value = 0.1422366992399565
$ws.Range("AQ5").Value = 248.1578947368421

# Row 3 no longer reports ebitdard/ebitda/operating/after-tax-operating margin
# or net margin (columns G, H, I, J, L) nor ebit_interest_expenses (AO)
$ws.Range("G3").ClearContents()
$ws.Range("H3").ClearContents()
$ws.Range("I3").ClearContents()
$ws.Range("J3").ClearContents()
$ws.Range("L3").ClearContents()
$ws.Range("AO3").ClearContents()

# The former row 6 (old "Weebit Nano Limited" entry) is dropped; Excel shifts
# nothing up underneath it and automatically shrinks the sheet dimension to AQ5
$ws.Rows.Item(6).Delete()
